$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (ChatBot): add Actual / Comment values
$ws.Range("E6").Value = "1 hour"
$ws.Range("F6").Value = "The chatbot answers research questions based on semantic search from uploaded journal chunks.  It shows the answer, top 5 citations with metadata and links, and a usage chart. The answer is LLM generated by understanding the context and rewriting in its own words."

# Row 9 (Generate summary): update comment text
$ws.Range("B9").Value = "Generate the summary of the user input"
